# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" worksheet (same shape as the existing
#    quarterly sheets) right before the "总计" (total) sheet.
# 2) Populate it with the two new fund holdings for 2022-Q1.
# 3) Update the "总计" sheet: prepend a 2022-Q1 summary row, pushing the
#    existing 2021-Q1 / 2020-Q4 rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet just before "总计"
# ---------------------------------------------------------------------
$q1src = $wb.Worksheets.Item("2021-Q1")

$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$newSheet.Name = "2022-Q1"

# NOTE: `Add(Before:=...)` splices the new sheet into the "总计" slot and
# pushes the original "总计" sheet one slot over, so the anchor captured
# in a variable *before* the call now resolves to the new sheet instead
# of the original one. Re-resolve "总计" by name afterwards to get a
# handle on the (now shifted) original sheet.
$totalSheet = $wb.Worksheets.Item("总计")

# Match the page setup used by the other quarterly sheets.
$newSheet.PageSetup.LeftMargin   = 54
$newSheet.PageSetup.RightMargin  = 54
$newSheet.PageSetup.TopMargin    = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36
$newSheet.Outline.SummaryRow    = 1
$newSheet.Outline.SummaryColumn = 1

# Carry over the header / index-column formatting from "2021-Q1".
$q1src.Range("B1:H1").Copy()
$newSheet.Range("B1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$q1src.Range("A2:A3").Copy()
$newSheet.Range("A2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Headers
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2 - 014887 招商安福1年定期开放债券
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'014887"
$newSheet.Range("B2").ClearFormats()
$newSheet.Range("C2").Value = "招商安福1年定期开放债券"
$newSheet.Range("D2").Value = "'17.22"
$newSheet.Range("D2").ClearFormats()
$newSheet.Range("E2").Value = "'27.65"
$newSheet.Range("E2").ClearFormats()
$newSheet.Range("F2").Value = "'2.15"
$newSheet.Range("F2").ClearFormats()
$newSheet.Range("G2").Value = "'0.3702"
$newSheet.Range("G2").ClearFormats()
$newSheet.Range("H2").Value = 3

# Row 3 - 005701 上投摩根香港精选港股通混合
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'005701"
$newSheet.Range("B3").ClearFormats()
$newSheet.Range("C3").Value = "上投摩根香港精选港股通混合"
$newSheet.Range("D3").Value = "'0.44"
$newSheet.Range("D3").ClearFormats()
$newSheet.Range("E3").Value = "'84.37"
$newSheet.Range("E3").ClearFormats()
$newSheet.Range("F3").Value = "'3.11"
$newSheet.Range("F3").ClearFormats()
$newSheet.Range("G3").Value = "'0.0137"
$newSheet.Range("G3").ClearFormats()
$newSheet.Range("H3").Value = 8

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet - add the 2022-Q1 row at the top, push the
#    rest down.
# ---------------------------------------------------------------------
$ws = $totalSheet

# Push 2021-Q1's row (currently row 2) down to row 4, and give it the
# same formatting as row 3 had originally.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "2020-Q4"
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 0.02

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "2021-Q1"
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 0.2

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q1"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 0.38

# Keep the originally-active first sheet selected (adding a sheet would
# otherwise leave the newly inserted sheet active/selected).
$wb.Worksheets.Item(1).Activate()

Write-Output "2022-Q1 sheet added and 总计 updated"
